$wb = $excel.ActiveWorkbook

# Delete the blank placeholder "Sheet1"
$wsDelete = $wb.Worksheets.Item("Sheet1")
$wsDelete.Delete() | Out-Null

# --- Accounts sheet: Bank Name / Current Balance / Asset Type ---
$wsAccounts = $wb.Worksheets.Item("Accounts")

$wsAccounts.Range("A1").Value = "Bank Name"
$wsAccounts.Range("B1").Value = "Current Balance"
$wsAccounts.Range("C1").Value = "Asset Type"

$wsAccounts.Range("A2").Value = "Deutsche Bank"
$wsAccounts.Range("B2").Value = 11409.3
$wsAccounts.Range("C2").Value = "Fixed"

$wsAccounts.Range("A3").Value = "Sparkasse"
$wsAccounts.Range("B3").Value = 7632.86
$wsAccounts.Range("C3").Value = "Liquid"

$wsAccounts.Range("A4").Value = "N26"
$wsAccounts.Range("B4").Value = 328.61
$wsAccounts.Range("C4").Value = "Liquid"

$wsAccounts.Range("A5").Value = "Cash"
$wsAccounts.Range("B5").Value = 69
$wsAccounts.Range("C5").Value = "Liquid"

# --- Wealth Allocation sheet: Class / Balance + Asset Type / Asset Total ---
$wsWealth = $wb.Worksheets.Item("Wealth Allocation")

$wsWealth.Range("A1").Value = "Class"
$wsWealth.Range("B1").Value = "Balance"
$wsWealth.Range("G1").Value = "Asset Type"
$wsWealth.Range("H1").Value = "Asset Total"

$wsWealth.Range("A2").Value = "Core Liquid"
$wsWealth.Range("B2").Value = 4000
$wsWealth.Range("G2").Value = "Liquid"
$wsWealth.Range("H2").Value = 8030.47

$wsWealth.Range("A3").Value = "Emergency Layer"
$wsWealth.Range("B3").Value = 3000
$wsWealth.Range("G3").Value = "Fixed"
$wsWealth.Range("H3").Value = 11409.3

$wsWealth.Range("A4").Value = "Rest Collector"
$wsWealth.Range("B4").Value = 1030.47

$wsWealth.Range("A5").Value = "Fixed Asset"
$wsWealth.Range("B5").Value = 11409.3

$wsWealth.Columns.Item(1).ColumnWidth = 14.796875

$wsWealth.Activate() | Out-Null
$wsWealth.Range("F6").Select() | Out-Null
